$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2-7 contain a date serial value (45180) that needs to be
# bumped by one day to 45181, keeping the existing date style/format intact.
for ($row = 2; $row -le 7; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
